$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Akbar"
$ws.Range("B6").Value = "Tester#123"
$ws.Range("C6").Value = 30

$ws.Range("C6").Select()
